$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a rolling weekly feed: the oldest weekly record (row 473)
# is dropped, existing rows 474-499 shift up by one, and a brand-new
# weekly record lands at row 499. Columns D,H,I,J,K,L,M,N,O,P,Q carry the
# per-record data; A,B,C,E,F,G,R are constant for this market/product and
# are left untouched.

# @{D=Fecha(serial); H=Variedad; I=Calidad; J=Volumen; K=PrecioMin; L=PrecioMax;
#   M=PrecioProm; N=UnidadComercializacion; O=Origen; P=PrecioPorKg; Q=KgOUnidades}
$filaDatos = @(
    @{ Fila=473; D=44699; H='Larga vida'; I='Primera'; J=300; K=11000; L=11000; M=11000; N='$/caja 10 kilos'; O='Región de Arica y Parinacota'; P=1100; Q=10 },
    @{ Fila=474; D=44636; H='Larga vida'; I='Primera'; J=400; K=10000; L=11000; M=10500; N='$/bandeja 18 kilos'; O='Región del Maule'; P=583; Q=18 },
    @{ Fila=475; D=44636; H='Larga vida'; I='Primera'; J=600; K=7000; L=7500; M=7250; N='$/caja 15 kilos'; O='Región del Maule'; P=483; Q=15 },
    @{ Fila=476; D=44636; H='Semiduro'; I='Primera'; J=400; K=11000; L=12000; M=11500; N='$/bandeja 18 kilos'; O='Región del Maule'; P=639; Q=18 },
    @{ Fila=477; D=44405; H='Larga vida'; I='Primera'; J=160; K=11500; L=12000; M=11750; N='$/bandeja 18 kilos'; O='Región de Arica y Parinacota'; P=653; Q=18 },
    @{ Fila=478; D=44405; H='Larga vida'; I='Primera'; J=120; K=6000; L=6500; M=6250; N='$/caja 10 kilos'; O='Región de Arica y Parinacota'; P=625; Q=10 },
    @{ Fila=479; D=44202; H='Larga vida'; I='Primera'; J=330; K=6000; L=7000; M=6455; N='$/caja 15 kilos'; O='Región del Maule'; P=430; Q=15 },
    @{ Fila=480; D=44273; H='Larga vida'; I='Primera'; J=160; K=4000; L=4500; M=4219; N='$/caja 15 kilos'; O='Región del Maule'; P=281; Q=15 },
    @{ Fila=481; D=44273; H='Larga vida'; I='Segunda'; J=150; K=3000; L=3500; M=3267; N='$/caja 15 kilos'; O='Región del Maule'; P=218; Q=15 },
    @{ Fila=482; D=44273; H='Semiduro'; I='Primera'; J=140; K=5500; L=6000; M=5786; N='$/bandeja 18 kilos'; O='Provincia de Diguillín'; P=321; Q=18 },
    @{ Fila=483; D=44273; H='Semiduro'; I='Segunda'; J=115; K=4500; L=5000; M=4717; N='$/bandeja 18 kilos'; O='Provincia de Diguillín'; P=262; Q=18 },
    @{ Fila=484; D=44777; H='Larga vida'; I='Primera'; J=300; K=8000; L=9000; M=8500; N='$/bandeja 18 kilos'; O='Región de Arica y Parinacota'; P=472; Q=18 },
    @{ Fila=485; D=44159; H='Larga vida'; I='Primera'; J=250; K=5500; L=6000; M=5740; N='$/caja 10 kilos'; O='Región de Arica y Parinacota'; P=574; Q=10 },
    @{ Fila=486; D=44159; H='Larga vida'; I='Primera'; J=250; K=8000; L=9000; M=8520; N='$/caja 15 kilos'; O='Provincia de Talca'; P=568; Q=15 },
    @{ Fila=487; D=44159; H='Larga vida'; I='Segunda'; J=90; K=4500; L=4500; M=4500; N='$/caja 10 kilos'; O='Región de Arica y Parinacota'; P=450; Q=10 },
    @{ Fila=488; D=44589; H='Larga vida'; I='Primera'; J=500; K=8000; L=8500; M=8250; N='$/bandeja 18 kilos'; O='Región del Maule'; P=458; Q=18 },
    @{ Fila=489; D=44589; H='Larga vida'; I='Primera'; J=1000; K=5000; L=5500; M=5250; N='$/caja 15 kilos'; O='Región del Maule'; P=350; Q=15 },
    @{ Fila=490; D=44263; H='Larga vida'; I='Primera'; J=250; K=8500; L=9000; M=8760; N='$/bandeja 18 kilos'; O='Provincia de Talca'; P=487; Q=18 },
    @{ Fila=491; D=44263; H='Semiduro'; I='Primera'; J=140; K=4000; L=4500; M=4214; N='$/caja 15 kilos'; O='Provincia de Diguillín'; P=281; Q=15 },
    @{ Fila=492; D=44263; H='Semiduro'; I='Segunda'; J=140; K=3000; L=3500; M=3286; N='$/caja 15 kilos'; O='Provincia de Diguillín'; P=219; Q=15 },
    @{ Fila=493; D=44263; H='Semiduro'; I='Tercera'; J=105; K=2000; L=2500; M=2286; N='$/caja 15 kilos'; O='Provincia de Diguillín'; P=152; Q=15 },
    @{ Fila=494; D=44309; H='Semiduro'; I='Primera'; J=160; K=4000; L=4500; M=4250; N='$/caja 15 kilos'; O='Provincia de Diguillín'; P=283; Q=15 },
    @{ Fila=495; D=44771; H='Larga vida'; I='Primera'; J=400; K=5000; L=5500; M=5250; N='$/caja 10 kilos'; O='Región de Arica y Parinacota'; P=525; Q=10 },
    @{ Fila=496; D=44267; H='Larga vida'; I='Primera'; J=210; K=6500; L=7000; M=6714; N='$/bandeja 18 kilos'; O='Provincia de Talca'; P=373; Q=18 },
    @{ Fila=497; D=44267; H='Semiduro'; I='Primera'; J=130; K=5500; L=6000; M=5808; N='$/bandeja 18 kilos'; O='Provincia de Diguillín'; P=323; Q=18 },
    @{ Fila=498; D=44267; H='Semiduro'; I='Segunda'; J=125; K=4500; L=5000; M=4740; N='$/bandeja 18 kilos'; O='Provincia de Diguillín'; P=263; Q=18 },
    @{ Fila=499; D=44826; H='Larga vida'; I='Primera'; J=600; K=6000; L=6500; M=6250; N='$/caja 10 kilos'; O='Región de Arica y Parinacota'; P=625; Q=10 }
)

foreach ($registro in $filaDatos) {
    $f = $registro.Fila
    $ws.Range("D" + $f).Value = $registro.D
    $ws.Range("H" + $f).Value = $registro.H
    $ws.Range("I" + $f).Value = $registro.I
    $ws.Range("J" + $f).Value = $registro.J
    $ws.Range("K" + $f).Value = $registro.K
    $ws.Range("L" + $f).Value = $registro.L
    $ws.Range("M" + $f).Value = $registro.M
    $ws.Range("N" + $f).Value = $registro.N
    $ws.Range("O" + $f).Value = $registro.O
    $ws.Range("P" + $f).Value = $registro.P
    $ws.Range("Q" + $f).Value = $registro.Q
}
